$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows with swapped coin identity/order ---
$ws.Range('B27').Value2 = 'Cosmos'
$ws.Range('C27').Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '11.66'
$ws.Range('E27').Value2 = '  -2.55%  '

$ws.Range('B28').Value2 = 'Dai'
$ws.Range('C28').Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '0.999'
$ws.Range('E28').Value2 = '  -1.61%  '

$ws.Range('B35').Value2 = 'Filecoin'
$ws.Range('C35').Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '5.76'
$ws.Range('E35').Value2 = '  -0.75%  '

$ws.Range('B36').Value2 = 'RenderToken'
$ws.Range('C36').Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '5.19'
$ws.Range('E36').Value2 = '  +12.32%  '

$ws.Range('B44').Value2 = 'Celestia'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '12.79'
$ws.Range('E44').Value2 = '  -6.71%  '

$ws.Range('B45').Value2 = 'FirstDigitalUSD'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '1.00'
$ws.Range('E45').Value2 = '  +0.16%  '

# --- Regular price/volume updates ---
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value2 = '42.634.11'
$ws.Range('E2').Value2 = '  -2.30%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value2 = '2.230.49'
$ws.Range('E3').Value2 = '  -2.06%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value2 = '1.01'
$ws.Range('E4').Value2 = '  +0.57%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '112.03'
$ws.Range('E5').Value2 = '  -7.27%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '296.48'
$ws.Range('E6').Value2 = '  +11.02%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '0.626'
$ws.Range('E7').Value2 = '  -3.66%  '

$ws.Range('E8').Value2 = '  +0.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.612'
$ws.Range('E9').Value2 = '  -2.66%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '44.69'
$ws.Range('E10').Value2 = '  -7.76%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '0.0920'
$ws.Range('E11').Value2 = '  -3.06%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '54.30'
$ws.Range('E12').Value2 = '  -0.19%  '

$ws.Range('E13').Value2 = '  -5.05%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '1.03'
$ws.Range('E14').Value2 = '  +12.12%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '0.104'
$ws.Range('E15').Value2 = '  -2.34%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '15.13'
$ws.Range('E16').Value2 = '  -3.57%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '2.562.81'
$ws.Range('E17').Value2 = '  -2.29%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '2.244.60'
$ws.Range('E18').Value2 = '  -1.20%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '42.583.28'
$ws.Range('E19').Value2 = '  -2.42%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '7.36'
$ws.Range('E20').Value2 = '  +5.93%  '

$ws.Range('E21').Value2 = '  -3.95%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '73.02'
$ws.Range('E22').Value2 = '  +0.90%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '3.50'
$ws.Range('E23').Value2 = '  +21.12%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '2.40'
$ws.Range('E24').Value2 = '  -0.65%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '230.26'
$ws.Range('E25').Value2 = '  -2.39%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '9.22'
$ws.Range('E26').Value2 = '  -3.14%  '

$ws.Range('E29').Value2 = '  -0.86%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '38.69'
$ws.Range('E30').Value2 = '  -10.76%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '3.26'
$ws.Range('E31').Value2 = '  -4.20%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '173.64'
$ws.Range('E32').Value2 = '  +0.12%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '21.11'
$ws.Range('E33').Value2 = '  -2.68%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '0.0900'
$ws.Range('E34').Value2 = '  -2.84%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '4.35'
$ws.Range('E37').Value2 = '  +1.71%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '0.126'
$ws.Range('E38').Value2 = '  -3.61%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '0.0379'
$ws.Range('E39').Value2 = '  -2.00%  '

$ws.Range('E40').Value2 = '  -4.76%  '

$ws.Range('E41').Value2 = '  -5.81%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '72.61'
$ws.Range('E42').Value2 = '  -1.98%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '0.236'
$ws.Range('E43').Value2 = '  -1.73%  '

$ws.Range('E46').Value2 = '  -4.68%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '5.45'
$ws.Range('E47').Value2 = '  -8.06%  '

$ws.Range('E48').Value2 = '  +4.01%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '103.51'
$ws.Range('E49').Value2 = '  +0.27%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '8.60'
$ws.Range('E50').Value2 = '  +1.25%  '

$ws.Range('E51').Value2 = '  +7.08%  '
